{"js": "// Generated edit script: replace arithmetic-expression text in each\n// table cell according to the old -> new text mapping scraped from the\n// canonical OOXML diff. The mapping is ordered exactly as the cells\n// appear in the document (row-major), but we match defensively by the\n// *current* cell text (consumed in order) so the script is resilient\n// to harmless reordering/whitespace differences.\nconst REPLACEMENTS = [[\"47+34=81\", \"90-30=60\"], [\"20-0=20\", \"56+38=94\"], [\"41+23=64\", \"16+21=37\"], [\"31-14=17\", \"31+27=58\"], [\"84-41=43\", \"52-20=32\"], [\"49+1=50\", \"77-16=61\"], [\"8+18=26\", \"77-22=55\"], [\"38+48=86\", \"16+37=53\"], [\"78-44=34\", \"26-5=21\"], [\"14+79=93\", \"56-14=42\"], [\"80-23=57\", \"38+39=77\"], [\"65-16=49\", \"44-28=16\"], [\"53-21=32\", \"92-67=25\"], [\"82-53=29\", \"88-50=38\"], [\"8+51=59\", \"95-65=30\"], [\"96-18=78\", \"9+85=94\"], [\"18+28=46\", \"27+16=43\"], [\"19+15=34\", \"11+75=86\"], [\"11+13=24\", \"2+11=13\"], [\"72+7=79\", \"47+31=78\"], [\"76-62=14\", \"2+75=77\"], [\"26+36=62\", \"84-58=26\"], [\"66-31=35\", \"68-41=27\"], [\"40+43=83\", \"31+27=58\"], [\"96-95=1\", \"27+71=98\"], [\"80-59=21\", \"5+20=25\"], [\"19+7=26\", \"5+64=69\"], [\"50-38=12\", \"98-69=29\"], [\"27-23=4\", \"43+38=81\"], [\"6+66=72\", \"23+61=84\"], [\"84-56=28\", \"43-0=43\"], [\"80-6=74\", \"13+47=60\"], [\"1+66=67\", \"58+25=83\"], [\"47-17=30\", \"51-48=3\"], [\"85-63=22\", \"34+55=89\"], [\"17-8=9\", \"47-31=16\"], [\"2+88=90\", \"36-33=3\"], [\"57-35=22\", \"13+24=37\"], [\"74+15=89\", \"41-1=40\"], [\"92-34=58\", \"84-46=38\"], [\"65+28=93\", \"12+11=23\"], [\"51+6=57\", \"58+34=92\"], [\"26+20=46\", \"10-9=1\"], [\"95+0=95\", \"12+75=87\"], [\"71-6=65\", \"81-62=19\"], [\"94-9=85\", \"90-11=79\"], [\"29-12=17\", \"43+36=79\"], [\"15+27=42\", \"94-21=73\"], [\"50-24=26\", \"25+50=75\"], [\"7+19=26\", \"34+18=52\"], [\"99-96=3\", \"50-3=47\"], [\"63-19=44\", \"58+7=65\"], [\"63+33=96\", \"1+60=61\"], [\"83-33=50\", \"4+33=37\"], [\"12+47=59\", \"60+10=70\"], [\"47+12=59\", \"84-37=47\"], [\"86-53=33\", \"99-19=80\"], [\"58+14=72\", \"92-22=70\"], [\"17+69=86\", \"15+20=35\"], [\"59+11=70\", \"23+50=73\"], [\"62+25=87\", \"60-52=8\"], [\"62+9=71\", \"8+60=68\"], [\"93-76=17\", \"28+62=90\"], [\"54-31=23\", \"47-39=8\"], [\"49-10=39\", \"1+44=45\"], [\"45+26=71\", \"61-55=6\"], [\"68+26=94\", \"10-5=5\"], [\"61-50=11\", \"9+20=29\"], [\"78-25=53\", \"82+12=94\"], [\"88-4=84\", \"52+3=55\"], [\"32+55=87\", \"14+7=21\"], [\"50+12=62\", \"64-14=50\"], [\"51-11=40\", \"79-36=43\"], [\"66+32=98\", \"26+39=65\"], [\"73+9=82\", \"73+11=84\"], [\"48+4=52\", \"11+41=52\"], [\"63-42=21\", \"5+42=47\"], [\"89-68=21\", \"24+0=24\"], [\"10+69=79\", \"34+7=41\"], [\"58+6=64\", \"49-48=1\"], [\"65+9=74\", \"32+2=34\"], [\"37+33=70\", \"59-40=19\"], [\"85-44=41\", \"64-11=53\"], [\"78-8=70\", \"35+14=49\"], [\"78+2=80\", \"82-1=81\"], [\"32+52=84\", \"11+63=74\"], [\"22+61=83\", \"94-48=46\"], [\"98-13=85\", \"42+35=77\"], [\"20+34=54\", \"77-32=45\"], [\"19-19=0\", \"46+14=60\"], [\"52-48=4\", \"92-72=20\"], [\"27+51=78\", \"41-22=19\"], [\"90-51=39\", \"87-47=40\"], [\"71+28=99\", \"18+9=27\"], [\"47+10=57\", \"35+36=71\"], [\"58+5=63\", \"75-22=53\"], [\"4+93=97\", \"35+10=45\"], [\"72-46=26\", \"9+50=59\"], [\"72-10=62\", \"69-51=18\"], [\"80-61=19\", \"94-62=32\"]];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length === 0) {\n  throw new Error(\"No table found in document body.\");\n}\n\nconst table = tables.items[0];\ntable.load(\"values,rowCount\");\nawait context.sync();\n\nconst oldValues = table.values;\nconst rowCount = oldValues.length;\nconst colCount = oldValues[0].length;\n\n// Build a queue per old-text so duplicate old strings (if any) are\n// consumed in the same left-to-right, top-to-bottom order they were\n// recorded in REPLACEMENTS.\nconst queues = new Map();\nfor (const [oldText, newText] of REPLACEMENTS) {\n  if (!queues.has(oldText)) queues.set(oldText, []);\n  queues.get(oldText).push(newText);\n}\n\nconst newValues = [];\nfor (let r = 0; r < rowCount; r++) {\n  const row = [];\n  for (let c = 0; c < colCount; c++) {\n    const cur = oldValues[r][c];\n    const q = queues.get(cur);\n    if (q && q.length > 0) {\n      row.push(q.shift());\n    } else {\n      row.push(cur);\n    }\n  }\n  newValues.push(row);\n}\n\ntable.values = newValues;\nawait context.sync();\n", "ps1": "# Generated edit script: replace arithmetic-expression text in each\n# table cell according to the old -> new text mapping scraped from the\n# canonical OOXML diff. All 100 \"old\" expressions are unique in this\n# document, so a simple text -> text lookup table is sufficient and\n# avoids relying on positional (row/column) assumptions.\n\n$map = @{}\n$map[\"47+34=81\"] = \"90-30=60\"\n$map[\"20-0=20\"] = \"56+38=94\"\n$map[\"41+23=64\"] = \"16+21=37\"\n$map[\"31-14=17\"] = \"31+27=58\"\n$map[\"84-41=43\"] = \"52-20=32\"\n$map[\"49+1=50\"] = \"77-16=61\"\n$map[\"8+18=26\"] = \"77-22=55\"\n$map[\"38+48=86\"] = \"16+37=53\"\n$map[\"78-44=34\"] = \"26-5=21\"\n$map[\"14+79=93\"] = \"56-14=42\"\n$map[\"80-23=57\"] = \"38+39=77\"\n$map[\"65-16=49\"] = \"44-28=16\"\n$map[\"53-21=32\"] = \"92-67=25\"\n$map[\"82-53=29\"] = \"88-50=38\"\n$map[\"8+51=59\"] = \"95-65=30\"\n$map[\"96-18=78\"] = \"9+85=94\"\n$map[\"18+28=46\"] = \"27+16=43\"\n$map[\"19+15=34\"] = \"11+75=86\"\n$map[\"11+13=24\"] = \"2+11=13\"\n$map[\"72+7=79\"] = \"47+31=78\"\n$map[\"76-62=14\"] = \"2+75=77\"\n$map[\"26+36=62\"] = \"84-58=26\"\n$map[\"66-31=35\"] = \"68-41=27\"\n$map[\"40+43=83\"] = \"31+27=58\"\n$map[\"96-95=1\"] = \"27+71=98\"\n$map[\"80-59=21\"] = \"5+20=25\"\n$map[\"19+7=26\"] = \"5+64=69\"\n$map[\"50-38=12\"] = \"98-69=29\"\n$map[\"27-23=4\"] = \"43+38=81\"\n$map[\"6+66=72\"] = \"23+61=84\"\n$map[\"84-56=28\"] = \"43-0=43\"\n$map[\"80-6=74\"] = \"13+47=60\"\n$map[\"1+66=67\"] = \"58+25=83\"\n$map[\"47-17=30\"] = \"51-48=3\"\n$map[\"85-63=22\"] = \"34+55=89\"\n$map[\"17-8=9\"] = \"47-31=16\"\n$map[\"2+88=90\"] = \"36-33=3\"\n$map[\"57-35=22\"] = \"13+24=37\"\n$map[\"74+15=89\"] = \"41-1=40\"\n$map[\"92-34=58\"] = \"84-46=38\"\n$map[\"65+28=93\"] = \"12+11=23\"\n$map[\"51+6=57\"] = \"58+34=92\"\n$map[\"26+20=46\"] = \"10-9=1\"\n$map[\"95+0=95\"] = \"12+75=87\"\n$map[\"71-6=65\"] = \"81-62=19\"\n$map[\"94-9=85\"] = \"90-11=79\"\n$map[\"29-12=17\"] = \"43+36=79\"\n$map[\"15+27=42\"] = \"94-21=73\"\n$map[\"50-24=26\"] = \"25+50=75\"\n$map[\"7+19=26\"] = \"34+18=52\"\n$map[\"99-96=3\"] = \"50-3=47\"\n$map[\"63-19=44\"] = \"58+7=65\"\n$map[\"63+33=96\"] = \"1+60=61\"\n$map[\"83-33=50\"] = \"4+33=37\"\n$map[\"12+47=59\"] = \"60+10=70\"\n$map[\"47+12=59\"] = \"84-37=47\"\n$map[\"86-53=33\"] = \"99-19=80\"\n$map[\"58+14=72\"] = \"92-22=70\"\n$map[\"17+69=86\"] = \"15+20=35\"\n$map[\"59+11=70\"] = \"23+50=73\"\n$map[\"62+25=87\"] = \"60-52=8\"\n$map[\"62+9=71\"] = \"8+60=68\"\n$map[\"93-76=17\"] = \"28+62=90\"\n$map[\"54-31=23\"] = \"47-39=8\"\n$map[\"49-10=39\"] = \"1+44=45\"\n$map[\"45+26=71\"] = \"61-55=6\"\n$map[\"68+26=94\"] = \"10-5=5\"\n$map[\"61-50=11\"] = \"9+20=29\"\n$map[\"78-25=53\"] = \"82+12=94\"\n$map[\"88-4=84\"] = \"52+3=55\"\n$map[\"32+55=87\"] = \"14+7=21\"\n$map[\"50+12=62\"] = \"64-14=50\"\n$map[\"51-11=40\"] = \"79-36=43\"\n$map[\"66+32=98\"] = \"26+39=65\"\n$map[\"73+9=82\"] = \"73+11=84\"\n$map[\"48+4=52\"] = \"11+41=52\"\n$map[\"63-42=21\"] = \"5+42=47\"\n$map[\"89-68=21\"] = \"24+0=24\"\n$map[\"10+69=79\"] = \"34+7=41\"\n$map[\"58+6=64\"] = \"49-48=1\"\n$map[\"65+9=74\"] = \"32+2=34\"\n$map[\"37+33=70\"] = \"59-40=19\"\n$map[\"85-44=41\"] = \"64-11=53\"\n$map[\"78-8=70\"] = \"35+14=49\"\n$map[\"78+2=80\"] = \"82-1=81\"\n$map[\"32+52=84\"] = \"11+63=74\"\n$map[\"22+61=83\"] = \"94-48=46\"\n$map[\"98-13=85\"] = \"42+35=77\"\n$map[\"20+34=54\"] = \"77-32=45\"\n$map[\"19-19=0\"] = \"46+14=60\"\n$map[\"52-48=4\"] = \"92-72=20\"\n$map[\"27+51=78\"] = \"41-22=19\"\n$map[\"90-51=39\"] = \"87-47=40\"\n$map[\"71+28=99\"] = \"18+9=27\"\n$map[\"47+10=57\"] = \"35+36=71\"\n$map[\"58+5=63\"] = \"75-22=53\"\n$map[\"4+93=97\"] = \"35+10=45\"\n$map[\"72-46=26\"] = \"9+50=59\"\n$map[\"72-10=62\"] = \"69-51=18\"\n$map[\"80-61=19\"] = \"94-62=32\"\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\nfor ($r = 1; $r -le $t.Rows.Count; $r++) {\n  for ($c = 1; $c -le $t.Columns.Count; $c++) {\n    $cell = $t.Cell($r, $c)\n    # Strip the trailing cell-end marker (CR + BEL) Word appends to cell text.\n    $cur = $cell.Range.Text -replace \"[\\r\\x07]+$\", \"\"\n    if ($map.ContainsKey($cur)) {\n      $cell.Range.Text = $map[$cur]\n    }\n  }\n}\n"}
